$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = -8.345299999999995
$ws.Range("D10").Value = -8.201899999999991
$ws.Range("D12").Value = -8.018600000000001
$ws.Range("D18").Value = -8.240099999999998
